$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value (to be broadcast across columns J:AS)
$updates = @{
    100 = 302691.0128
    101 = 70373.84845
    102 = 242525.1209
    103 = 166959.6981
    104 = 232839.2406
    105 = 444122.8059
    106 = 26677.7633
    107 = 518706.2888
    114 = 2713.795026
    115 = 689759.2613
}

foreach ($row in $updates.Keys) {
    $value = $updates[$row]
    $range = $ws.Range("J$row`:AS$row")
    $range.Value = $value
}
